$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.429.80"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").Value = "2.519.44"
$ws.Range("E3").Value = "  +1.09%  "

$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.62"
$ws.Range("E5").Value = "  -0.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.63"
$ws.Range("E6").Value = "  -0.92%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("E8").Value = "  -1.03%  "

$ws.Range("D9").Value = "2.542.59"
$ws.Range("E9").Value = "  +0.86%  "

$ws.Range("E10").Value = "  -0.12%  "

$ws.Range("E11").Value = "  +0.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.60"
$ws.Range("E12").Value = "  +0.54%  "

$ws.Range("E13").Value = "  +2.47%  "

$ws.Range("D14").Value = "2.965.44"
$ws.Range("E14").Value = "  +0.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.69"
$ws.Range("E15").Value = "  -3.76%  "

$ws.Range("D16").Value = "59.325.72"
$ws.Range("E16").Value = "  -0.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000142"
$ws.Range("E17").Value = "  +0.94%  "

$ws.Range("D18").Value = "2.531.25"
$ws.Range("E18").Value = "  +0.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.21"
$ws.Range("E19").Value = "  -1.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.29"
$ws.Range("E20").Value = "  -1.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.76"
$ws.Range("E21").Value = "  -0.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +1.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.88"
$ws.Range("E23").Value = "  +0.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.13"
$ws.Range("E24").Value = "  +1.15%  "

$ws.Range("E25").Value = "  -3.92%  "

$ws.Range("E26").Value = "  +1.94%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.994"
$ws.Range("E27").Value = "  -0.93%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.18"
$ws.Range("E28").Value = "  +4.80%  "

$ws.Range("D29").Value = "0.0₃0788"
$ws.Range("E29").Value = "  -1.51%  "

$ws.Range("E30").Value = "  -0.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.73"
$ws.Range("E31").Value = "  -0.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.21"
$ws.Range("E32").Value = "  -6.23%  "

$ws.Range("E33").Value = "  +3.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.00"
$ws.Range("E35").Value = "  +0.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.76"
$ws.Range("E36").Value = "  -0.83%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.41"
$ws.Range("E37").Value = "  -2.57%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.61"
$ws.Range("E38").Value = "  -7.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.00"
$ws.Range("E39").Value = "  +0.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.60"
$ws.Range("E40").Value = "  -4.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.837"
$ws.Range("E41").Value = "  -0.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.71"
$ws.Range("E42").Value = "  -2.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "295.60"
$ws.Range("E43").Value = "  -7.18%  "

$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("E45").Value = "  -0.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.82"
$ws.Range("E46").Value = "  +1.25%  "

$ws.Range("E47").Value = "  -0.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.86"
$ws.Range("E48").Value = "  +0.68%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.26"
$ws.Range("E49").Value = "  -2.64%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0228"
$ws.Range("E50").Value = "  -2.27%  "

$ws.Range("B51").Value = "Hedera"
$ws.Range("C51").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0515"
$ws.Range("E51").Value = "  -3.63%  "
